$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the duplicate row 13 (a partial copy of row 14) entirely,
# shifting all rows below it up by one.
$ws.Rows.Item(13).Delete()

# Update the selection to reflect where the user clicked after the edit.
$ws.Range("A8").Select()
